$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.5141455444626758, 0.9897622703593392)
    3  = @(0.2419097835045483, 0.995271867197829)
    4  = @(0.338055960087607, 0.9934965550431772)
    5  = @(0.4153349664147831, 0.9918101631065696)
    6  = @(0.7292360521529262, 0.9785804268676348)
    7  = @(0.404058060215838, 0.9944158082032909)
    8  = @(0.1526505001906213, 0.9984266638436701)
    9  = @(0.5254954434052722, 0.9968697735794284)
    10 = @(0.1045862308324393, 0.9980894572443918)
    11 = @(0.3012841367677325, 0.9777400304056636)
    12 = @(0.04999187914524317, 0.9985298241229227)
    13 = @(0.08414612021499523, 0.9992018363398451)
    14 = @(0.08745618430458782, 0.9988158972495058)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}
